# Apply line-wrapping / whitespace-cleanup edits to course_catalogue sheet
# per commit "stler & linter on scripts" - reflows several long description
# cells onto multiple lines and normalises a couple of name-list separators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Introduction to Natural Language Programming with Python
$ws.Range("D7").Value = "Natural Language Processing is a sub-field of Artificial Intelligence.`nIt is used for processing and analysing large amounts of natural language.`nSome applications include search engines (Google), text classification (spam filters),`nidentifying sentiments for a product (sentiment analysis), methods for discovering`nabstract topics in a collection of documents (topic modelling) and machine translation technologies.`nConcepts covered include cleaning, exploring datasets through methods rooted in Corpus Linguistics,`nand application of feature engineering techniques to transform textual data into a numerical representation.`nKey techniques such as word embeddings and language modelling are also introduced as well as illustrations`nas to how they can be performed over a dataset."

# Row 8 - Introduction to Python Programming
$ws.Range("D8").Value = "This course is delivered through the Jupyter notebook application.`nIt begins with a coverage of fundamental building blocks in Python -  numeric data types, strings, lists,`ndictionaries, sets - replete with examples. Illustrations are then provided on the use of these data types`nto compose code with selection and iteration constructs.  To promote modular and readable code the set-up`nand use of functions with parameters are also covered."

$ws.Range("E8").Value = "Participants should attain a good understanding of basic data types in Python and associated methods`nand constructs that can be applied to them."

$ws.Range("K8").Value = "Richard Leyshon,`nKaveh Jahanshahi"

# Row 9 - Fundemental theories in Machine Learning
$ws.Range("K9").Value = "Laurie Baker `nIsabela Breton"

# Row 10 - Natural Language Processing with R
$ws.Range("D10").Value = "Natural Language Processing is a sub-field of Artificial Intelligence.`nIt is used for processing and analysing large amounts of natural language (texts).`nSome applications include search engines (Google), text classification (spam filters),`nidentifying sentiments for a product (sentiment analysis), methods for discovering`nabstract topics in a collection of documents (topic modelling) and machine translation technologies.`nIn this course you learn about exploratory analysis of text data, introduced to sentiment analysis of texts using sentiment lexicons and the concept of topic modelling (package topicmodels)."

# Row 11 - Introduction to Reproducibility
$ws.Range("J11").Value = "Laurie L. Baker`nRichard Leyshon"
